$wb = $excel.ActiveWorkbook

# --- Sheet: Forecast Comparison ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

$ws1.Range("L2").Value = 0.87

$ws1.Range("D3").Value = 0
$ws1.Range("H3").Value = 85.76000000000001
$ws1.Range("L3").Value = 1.04

$ws1.Range("H4").Value = 70.63
$ws1.Range("L4").Value = 1.05

$ws1.Range("H5").Value = 69.63
$ws1.Range("L5").Value = 1.04

$ws1.Range("H6").Value = 68.63
$ws1.Range("L6").Value = 1.03

$ws1.Range("H7").Value = 67.63
$ws1.Range("L7").Value = 1.2

$ws1.Range("H8").Value = 53.7
$ws1.Range("L8").Value = 0.99

$ws1.Range("H9").Value = 52.7
$ws1.Range("L9").Value = 1.17

$ws1.Range("H10").Value = 64.15000000000001
$ws1.Range("L10").Value = 1.17

$ws1.Range("H11").Value = 62
$ws1.Range("L11").Value = 0.96

$ws1.Range("H12").Value = 61
$ws1.Range("L12").Value = 1.06

$ws1.Range("H13").Value = 60
$ws1.Range("L13").Value = 1.09

$ws1.Range("H14").Value = 59
$ws1.Range("L14").Value = 0.99

$ws1.Range("D15").Value = 0
$ws1.Range("H15").Value = 70.89
$ws1.Range("L15").Value = 0.96

$ws1.Range("D16").Value = 0
$ws1.Range("H16").Value = 69.89
$ws1.Range("L16").Value = 0.91

$ws1.Range("D17").Value = 0
$ws1.Range("H17").Value = 86.11
$ws1.Range("L17").Value = 1.2

# --- Sheet: Summary ---
$ws2 = $wb.Worksheets.Item("Summary")
$ws2.Range("B9").Value = "17"
